$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -181022.72
$ws.Range("B3").Value = 338613.39
$ws.Range("B4").Value = -179006
$ws.Range("B5").Value = 2254401.67
$ws.Range("B6").Value = -80719.32000000001
$ws.Range("B7").Value = -190509.04
$ws.Range("B8").Value = -80232.52
$ws.Range("B9").Value = -279746.85
$ws.Range("B10").Value = -141099.22
$ws.Range("B11").Value = 965763.85
$ws.Range("B12").Value = -965763.85
$ws.Range("B13").Value = -83002.44
$ws.Range("B14").Value = 413324.54
$ws.Range("B15").Value = 1791001.49
